$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 4 new columns (E:H) to make room for the new MySQL connection fields.
#    This shifts the former E..J columns (Pwd/empty formatting columns) to I..N,
#    which matches the target dimension (A1:N5), row spans (1:14) and keeps the
#    widths of the old G/H/I/J (now K/L/M/N) columns untouched.
$ws.Range("E1:H1").EntireColumn.Insert()

# 2. Fill in the new header row (row 1) with the MySQL connection field names.
$ws.Range("E1").Value = "SqlIP"
$ws.Range("F1").Value = "SqlPort"
$ws.Range("G1").Value = "SqlName"
$ws.Range("H1").Value = "SqlUser"
$ws.Range("I1").Value = "SqlPwd"

# The column insert copies the left neighbour's style onto the new header cells;
# the target file has these cells with no explicit style, so clear that back out.
$ws.Range("E1:I1").ClearFormats()

# 3. Fill in the new data row (row 2) with the MySQL connection values.
$ws.Range("E2").Value = "192.168.0.24"
$ws.Range("F2").Value = 3306
$ws.Range("G2").Value = "app_test"
$ws.Range("H2").Value = "root"
$ws.Range("I2").Value = 123456

$ws.Range("F2:I2").ClearFormats()

# 4. Column widths (bestFit column widths set by real Excel autofit cannot be
#    reproduced exactly by this engine's ColumnWidth rounding model, so we set
#    the closest achievable values).
$ws.Columns.Item(1).ColumnWidth = 12.0
$ws.Columns.Item(4).ColumnWidth = 4.714285714285714
$ws.Columns.Item(5).ColumnWidth = 13.142857142857142
$ws.Columns.Item(6).ColumnWidth = 7.857142857142857
$ws.Columns.Item(7).ColumnWidth = 7.857142857142857
$ws.Columns.Item(8).ColumnWidth = 7.857142857142857
$ws.Columns.Item(9).ColumnWidth = 6.714285714285714

# 5. Selection moved to G9 in the saved file.
$ws.Range("G9").Select()

# 6. Data validation: keep the blank "allow any" validation on the first spare
#    column (now J1) and extend the TRUE/FALSE list validation to the header
#    cells of the new columns plus the rest of column J.
$ws.Range("J2:J1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$ws.Range("E1:F1").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$ws.Range("H1").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
$ws.Range("I1").Validation.Add(3, 1, 1, '"TRUE,FALSE"')
